# Ultimate fix on backend
# Swap the "tp" (column C) and "fp" (column D) values for each data row,
# then recompute the "precision" (column G) and "fscore" (column I)
# columns based on the corrected tp/fp values. "recall" (column H) is
# unaffected because "fn" (column F) is always 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $tpCell = $ws.Cells.Item($r, 3)  # column C, tp
    $fpCell = $ws.Cells.Item($r, 4)  # column D, fp
    $fnCell = $ws.Cells.Item($r, 6)  # column F, fn

    $oldTp = $tpCell.Value2
    $oldFp = $fpCell.Value2
    $fn = $fnCell.Value2

    # swap tp and fp
    $newTp = $oldFp
    $newFp = $oldTp

    $tpCell.Value = $newTp
    $fpCell.Value = $newFp

    $precision = $newTp / ($newTp + $newFp)
    if (($newTp + $fn) -ne 0) {
        $recall = $newTp / ($newTp + $fn)
    } else {
        $recall = 1
    }
    $fscore = 2 * $precision * $recall / ($precision + $recall)

    $ws.Cells.Item($r, 7).Value = $precision   # column G, precision
    $ws.Cells.Item($r, 9).Value = $fscore      # column I, fscore
}
